$wb = $excel.ActiveWorkbook

# =====================================================================
# Sheet 1: Overview - update row 2, append row 3 for the new file entry
# =====================================================================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md'
$ws1.Range("C2").Value = '.md'
$ws1.Range("E2").Value = 'Handed back: in sync with en-US'
$ws1.Range("F2").Value = 'Handed back: in sync with en-US'
$ws1.Range("G2").Value = '2016-08-24 17:04:13'

$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

$ws1.Range("A3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.md'
$ws1.Range("C3").Value = '.md'
$ws1.Range("E3").Value = 'Handed back: in sync with en-US'
$ws1.Range("F3").Value = 'Handed back: in sync with en-US'
$ws1.Range("G3").Value = '2016-08-24 17:04:13'

# Rebuild the hyperlinks on the Overview sheet (B2, B3)
$ws1.Range("A1").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2112119329fb84284b499b85e513faaad98c0ab1/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md', "", "", 'e2e\383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md') | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2112119329fb84284b499b85e513faaad98c0ab1/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md', "", "", 'e2e\42647d57-8228-4722-a6e3-4fd76a0d03a6.md') | Out-Null

# =====================================================================
# Sheet 2: zh-cn - update row 2, append row 3 for the new file entry
# =====================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md'
$ws2.Range("B2").Value = '.md'
$ws2.Range("C2").Value = 'Handed back: in sync with en-US'
$ws2.Range("D2").Value = 'e2e'
$ws2.Range("E2").Value = 'ht'
$ws2.Range("F2").Value = 'False'
$ws2.Range("G2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.9c652668270faba77ec5a33cb84f14b0ab7c5182.zh-cn.xlf'
$ws2.Range("H2").Value = '2016-08-24 17:04:01'
$ws2.Range("I2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md'
$ws2.Range("J2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.9c652668270faba77ec5a33cb84f14b0ab7c5182.zh-cn.xlf'
$ws2.Range("K2").Value = '2016-08-24 17:04:30'
$ws2.Range("M2").Value = 'True'
$ws2.Range("O2").Value = 'False'

$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

$ws2.Range("A3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.md'
$ws2.Range("B3").Value = '.md'
$ws2.Range("C3").Value = 'Handed back: in sync with en-US'
$ws2.Range("D3").Value = 'e2e'
$ws2.Range("E3").Value = 'ht'
$ws2.Range("F3").Value = 'True'
$ws2.Range("G3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.8f99d4fb6fa6ab4fea68a2aee59820cebba79779.zh-cn.xlf'
$ws2.Range("H3").Value = '2016-08-24 17:04:01'
$ws2.Range("I3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.md'
$ws2.Range("J3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.8f99d4fb6fa6ab4fea68a2aee59820cebba79779.zh-cn.xlf'
$ws2.Range("K3").Value = '2016-08-24 17:04:30'
$ws2.Range("M3").Value = 'True'
$ws2.Range("O3").Value = 'False'

# Blank cells (L, N, P) still need to physically exist on both rows
$ws2.Range("L2").NumberFormat = "General"
$ws2.Range("N2").NumberFormat = "General"
$ws2.Range("P2").NumberFormat = "General"
$ws2.Range("L3").NumberFormat = "General"
$ws2.Range("N3").NumberFormat = "General"
$ws2.Range("P3").NumberFormat = "General"

# Apply date/time display format to the datetime columns (H, K)
$ws2.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild the hyperlinks on the zh-cn sheet (A2, I2, A3, I3)
$ws2.Range("A1").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2112119329fb84284b499b85e513faaad98c0ab1/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md', "", "", '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md') | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/07caaa89e899fcb7eb3e9d37255ae372be6c8e8d/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md', "", "", '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md') | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2112119329fb84284b499b85e513faaad98c0ab1/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md', "", "", '42647d57-8228-4722-a6e3-4fd76a0d03a6.md') | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/07caaa89e899fcb7eb3e9d37255ae372be6c8e8d/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md', "", "", '42647d57-8228-4722-a6e3-4fd76a0d03a6.md') | Out-Null

# =====================================================================
# Sheet 3: de-de - update row 2, append row 3 for the new file entry
# =====================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md'
$ws3.Range("B2").Value = '.md'
$ws3.Range("C2").Value = 'Handed back: in sync with en-US'
$ws3.Range("D2").Value = 'e2e'
$ws3.Range("E2").Value = 'ht'
$ws3.Range("F2").Value = 'False'
$ws3.Range("G2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.9c652668270faba77ec5a33cb84f14b0ab7c5182.de-de.xlf'
$ws3.Range("H2").Value = '2016-08-24 17:04:13'
$ws3.Range("I2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md'
$ws3.Range("J2").Value = '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.9c652668270faba77ec5a33cb84f14b0ab7c5182.de-de.xlf'
$ws3.Range("K2").Value = '2016-08-24 17:04:40'
$ws3.Range("M2").Value = 'True'
$ws3.Range("O2").Value = 'False'

$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

$ws3.Range("A3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.md'
$ws3.Range("B3").Value = '.md'
$ws3.Range("C3").Value = 'Handed back: in sync with en-US'
$ws3.Range("D3").Value = 'e2e'
$ws3.Range("E3").Value = 'ht'
$ws3.Range("F3").Value = 'True'
$ws3.Range("G3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.8f99d4fb6fa6ab4fea68a2aee59820cebba79779.de-de.xlf'
$ws3.Range("H3").Value = '2016-08-24 17:04:13'
$ws3.Range("I3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.md'
$ws3.Range("J3").Value = '42647d57-8228-4722-a6e3-4fd76a0d03a6.8f99d4fb6fa6ab4fea68a2aee59820cebba79779.de-de.xlf'
$ws3.Range("K3").Value = '2016-08-24 17:04:40'
$ws3.Range("M3").Value = 'True'
$ws3.Range("O3").Value = 'False'

# Blank cells (L, N, P) still need to physically exist on both rows
$ws3.Range("L2").NumberFormat = "General"
$ws3.Range("N2").NumberFormat = "General"
$ws3.Range("P2").NumberFormat = "General"
$ws3.Range("L3").NumberFormat = "General"
$ws3.Range("N3").NumberFormat = "General"
$ws3.Range("P3").NumberFormat = "General"

# Apply date/time display format to the datetime columns (H, K)
$ws3.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild the hyperlinks on the de-de sheet (A2, I2, A3, I3)
$ws3.Range("A1").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2112119329fb84284b499b85e513faaad98c0ab1/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md', "", "", '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md') | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/54a036193cd78dc595372c62115a6cdf13387a03/e2e/383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md', "", "", '383509dd-8e05-4b7e-8273-1d9fa6c3dfe6.md') | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2112119329fb84284b499b85e513faaad98c0ab1/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md', "", "", '42647d57-8228-4722-a6e3-4fd76a0d03a6.md') | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/54a036193cd78dc595372c62115a6cdf13387a03/e2e/42647d57-8228-4722-a6e3-4fd76a0d03a6.md', "", "", '42647d57-8228-4722-a6e3-4fd76a0d03a6.md') | Out-Null

